$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update URL
$ws.Range("B2").Value = "https://prod-unilevel.epixel.link/en/register/"

# Update Sponsor
$ws.Range("B3").Value = "mpfp-base-unilevel-business-admin,user1,user2"

# Update Subdomain and add a note in C6
$ws.Range("B6").Value = "antp087123"
$ws.Range("C6").Value = "If needed change"

# Update Enrollment Package
$ws.Range("B11").Value = "Bronze"

# Move the blank marker row from row 19 to row 23
$ws.Range("A19:AA19").ClearContents()
$ws.Range("B23:AA23").Value = ""
